# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> new F-column value (identical update set for both data sheets)
$updates = @{
    3  = 1401
    5  = 260
    8  = 18
    9  = 190
    11 = 4697
    12 = 6956
    18 = 4165
    19 = 904
    22 = 2747
    24 = 552
    26 = 383
    27 = 381
    29 = 241
    30 = 49
    32 = 1049
    34 = 522
    37 = 5
    41 = 213
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
